$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 16 and 17: copy the formatting of row 15 first, then fill values ---
$ws.Range("A15:E15").Copy() | Out-Null
$ws.Range("A16:E16").PasteSpecial(-4122) | Out-Null
$ws.Range("A17:E17").PasteSpecial(-4122) | Out-Null

# --- Rows 8-15: names shift (sharedStrings gained "line7"/"line8"), and the
#     from_bus/to_bus/in_service results were recalculated ---
$ws.Cells.Item(8, 2).Value  = "line7"
$ws.Cells.Item(8, 3).Value  = 14
$ws.Cells.Item(8, 4).Value  = 11

$ws.Cells.Item(9, 2).Value  = "line8"
$ws.Cells.Item(9, 3).Value  = 16

$ws.Cells.Item(10, 2).Value = "extr1"
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12

$ws.Cells.Item(11, 2).Value = "extr2"
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9

$ws.Cells.Item(12, 2).Value = "extr3"
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 5).Value = $True

$ws.Cells.Item(13, 2).Value = "extr4"
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $False

$ws.Cells.Item(14, 2).Value = "extr5"
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $False

$ws.Cells.Item(15, 2).Value = "extr6"
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11

# --- New rows 16 and 17 values ---
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "extr7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $False

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "extr8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $False
